$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column writes stay text (matches source formatting which mixes
# thousand-dot-separated and decimal-only strings) instead of Excel auto-numeric coercion.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '68.705.53'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '3.865.96'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '602.69'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '168.79'
$ws.Range('E6').Value = '  +2.84%  '
$ws.Range('D7').Value = '3.866.06'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').Value = '6.37'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('D14').Value = '37.69'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '4.513.69'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').Value = '3.870.42'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '68.813.63'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').Value = '18.43'
$ws.Range('E19').Value = '  +7.46%  '
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').Value = '10.82'
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('D22').Value = '479.16'
$ws.Range('E22').Value = '  -1.93%  '
$ws.Range('D23').Value = '0.741'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('E24').Value = '  -2.02%  '
$ws.Range('D25').Value = '84.90'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').Value = '2.26'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Value = '12.39'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('D28').Value = '10.10'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = '4.015.62'
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('D32').Value = '7.79'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').Value = '2.33'
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('D34').Value = '31.27'
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('D35').Value = '3.832.48'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '3.38'
$ws.Range('E37').Value = '  +11.14%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = '6.03'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.141'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').Value = '1.02'
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').Value = '0.319'
$ws.Range('E42').Value = '  -0.65%  '
$ws.Range('D43').Value = '2.03'
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').Value = '428.81'
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').Value = '47.82'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D47').Value = '8.65'
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('E48').Value = '  +14.50%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0362'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '142.41'
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('D51').Value = '39.32'
$ws.Range('E51').Value = '  +0.93%  '

$ws.Range("D2:D51").Style = "Normal"
